$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Slit3", "Robo2", "ECs", 2, 0.6666666666666666, 2.096169666666666, 6.288508999999999, 0.01755614486036394, 0.01755614486036394, 2, 0.6666666666666666, 1.847863, 5.543589, 0.7296835394802672, 0.7296835394802672, 3.873434368755666, 34.860909318801, 0.01281042992133866, 0.01281042992133866),
    @("ECs", "Slit3", "Robo2", "FAPs", 2, 0.6666666666666666, 2.096169666666666, 6.288508999999999, 0.01755614486036394, 0.01755614486036394, 3, 1, 0.661303, 1.983909, 0.261135113213977, 0.261135113213977, 1.386203289075667, 12.475829601681, 0.004584525875712117, 0.004584525875712117),
    @("ECs", "Slit3", "Robo2", "sCs", 2, 0.6666666666666666, 2.096169666666666, 6.288508999999999, 0.01755614486036394, 0.01755614486036394, 2, 0.6666666666666666, 0.023251, 0.069753, 0.009181347305755726, 0.009181347305755726, 0.04873804091966666, 0.4386423682769999, 0.0001611890633131597, 0.0001611890633131597),
    @("FAPs", "Slit3", "Robo2", "ECs", 3, 1, 108.7065656666667, 326.119697, 0.9104550287436967, 0.9104550287436968, 2, 0.6666666666666666, 1.847863, 5.543589, 0.7296835394802672, 0.7296835394802672, 200.8748405525037, 1807.873564972533, 0.664344047911309, 0.6643440479113091),
    @("FAPs", "Slit3", "Robo2", "FAPs", 3, 1, 108.7065656666667, 326.119697, 0.9104550287436967, 0.9104550287436968, 3, 1, 0.661303, 1.983909, 0.261135113213977, 0.261135113213977, 71.88797799506366, 646.9918019555729, 0.2377517770072199, 0.2377517770072199),
    @("FAPs", "Slit3", "Robo2", "sCs", 3, 1, 108.7065656666667, 326.119697, 0.9104550287436967, 0.9104550287436968, 2, 0.6666666666666666, 0.023251, 0.069753, 0.009181347305755726, 0.009181347305755726, 2.527536358315666, 22.747827224841, 0.008359203825167692, 0.008359203825167692),
    @("sCs", "Slit3", "Robo2", "ECs", 3, 1, 8.595326333333334, 25.785979, 0.07198882639593941, 0.07198882639593941, 2, 0.6666666666666666, 1.847863, 5.543589, 0.7296835394802672, 0.7296835394802672, 15.88298550429234, 142.946869538631, 0.05252906164761956, 0.05252906164761956),
    @("sCs", "Slit3", "Robo2", "FAPs", 3, 1, 8.595326333333334, 25.785979, 0.07198882639593941, 0.07198882639593941, 3, 1, 0.661303, 1.983909, 0.261135113213977, 0.261135113213977, 5.684115090212334, 51.157035811911, 0.01879881033104497, 0.01879881033104497),
    @("sCs", "Slit3", "Robo2", "sCs", 3, 1, 8.595326333333334, 25.785979, 0.07198882639593941, 0.07198882639593941, 2, 0.6666666666666666, 0.023251, 0.069753, 0.009181347305755726, 0.009181347305755726, 0.1998499325763333, 1.798649393187, 0.000660954417274875, 0.000660954417274875)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    for ($j = 0; $j -lt $data[$i].Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $data[$i][$j]
    }
}
